# Generate Report for Handoff
# Update the "Latest Handoff Datetime" for the 8c298c2c-... record (row 5)
# on both the zh-cn and de-de localization status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-19 04:33:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-19 04:33:27"
